$wb = $excel.ActiveWorkbook

# --- Fix testdata bug: posts that reference the bogus "null" space should
#     have their space column cleared instead of literally saying "null" ---
$posts = $wb.Worksheets.Item("posts")
$nullRows = @(2, 4, 5, 9, 12, 15, 19, 24, 25, 27)
foreach ($r in $nullRows) {
    $posts.Cells.Item($r, 4).Value = $null
}

# --- Updated post and repost templates: move the active tab/selection
#     from "users" to "posts" ---
$posts.Activate()
$posts.Range("D27").Select()
